$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None`n`nMSG: No decision has been made about which movie to show on Friday.`n"
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding Friday's movie.`n"
$ws.Range("C4").Value = "MSG: None`n`nMSG: The rights for `"Oppenheimer`" have been successfully acquired.`n"
$ws.Range("C5").Value = "MSG: None`n`nMSG: The decision has been recorded as no consensus was reached regarding a movie to show on Friday.`n"
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Oppenheimer`" has been recorded successfully.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been made.`n"
$ws.Range("D7").Value = "Barbie_was_selected, "
$ws.Range("C9").Value = "MSG: None`n`nMSG: The function for no decision has been called. No movie has been selected for Friday.`n"
$ws.Range("C10").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie.`"`n"
$ws.Range("C11").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has resulted in no conclusion being reached.`n"
$ws.Range("C12").Value = "MSG: None`n`nMSG: The decision has been recorded, and `"Oppenheimer`" will be shown on Friday.`n"
$ws.Range("C13").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been successfully recorded.`n"
$ws.Range("C14").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision was made regarding which movie to show on Friday.`n"
$ws.Range("C15").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie`" to be shown on Friday.`n"
$ws.Range("C16").Value = "MSG: None`n`nMSG: The decision has been recorded as no agreement was reached regarding the movie to be shown on Friday.`n"
$ws.Range("C17").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected during the committee's discussion.`n"
$ws.Range("C18").Value = "MSG: None`n`nMSG: The decision has been recorded with no movie selected for Friday.`n"
$ws.Range("C19").Value = "MSG: None`n`nMSG: The decision has been recorded with no movie selected for Friday.`n"
$ws.Range("C20").Value = "MSG: None`n`nMSG: The decision-making process did not result in a selection for Friday's movie.`n"
$ws.Range("C21").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no consensus was reached on the movie selection for Friday.`n"
$ws.Range("C22").Value = "MSG: None`n`nMSG: The decision has been recorded and the rights for `"Barbie`" have been acquired.`n"
$ws.Range("C23").Value = "MSG: None`n`nMSG: The decision has been made with no clear choice of a movie to show on Friday.`n"
$ws.Range("C24").Value = "MSG: None`n`nMSG: The decision-making process did not result in an agreement about which movie to show on Friday. Therefore, the outcome is classified as a `"no decision.`"`n"
$ws.Range("D24").Value = "no_decision, "
$ws.Range("C25").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision.`n"
$ws.Range("C26").Value = "MSG: None`n`nMSG: The decision resulted in no agreement on a movie for Friday.`n"
$ws.Range("C27").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been confirmed.`n"
$ws.Range("C28").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to show `"Barbie`" on Friday.`n"
$ws.Range("C29").Value = "MSG: None`n`nMSG: The call to the no_decision function has been successfully made, indicating that no decision was reached regarding the movie to be shown on Friday.`n"
$ws.Range("C30").Value = "MSG: None`n`nMSG: No decision was made regarding the movie for Friday.`n"
$ws.Range("C31").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie.`"`n"
$ws.Range("C32").Value = "MSG: None`n`nMSG: The decision process concluded without an agreement on the movie for Friday, so no acquisition is made.`n"
$ws.Range("C33").Value = "MSG: None`n`nMSG: The committee did not reach a decision regarding which movie to show on Friday.`n"
$ws.Range("C34").Value = "MSG: None`n`nMSG: No decision was made regarding which movie to show on Friday.`n"
$ws.Range("C35").Value = "MSG: None`n`nMSG: The rights to both movies have been acquired.`n"
$ws.Range("D35").Value = "both_movies, "
$ws.Range("C36").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights to `"Barbie.`"`n"
$ws.Range("C37").Value = "MSG: None`n`nMSG: The decision resulted in no movie being selected for Friday. If you have any further questions or need assistance, feel free to ask!`n"
$ws.Range("C38").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie.`"`n"
$ws.Range("D38").Value = "Barbie_was_selected, "
$ws.Range("C39").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire rights for the movie `"Barbie.`"`n"
$ws.Range("C40").Value = "MSG: None`n`nMSG: The decision process ended without a clear agreement on which movie to show on Friday.`n"
$ws.Range("C41").Value = "MSG: None`n`nMSG: The decision has been recorded, indicating that no movie was selected for Friday.`n"
$ws.Range("C42").Value = "MSG: None`n`nMSG: The decision regarding the movie for Friday was not made, so no action is recorded.`n"
$ws.Range("C43").Value = "MSG: None`n`nMSG: The decision has been recorded, and the movie `"Barbie`" will be acquired for Friday's showing.`n"
$ws.Range("C44").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday resulted in no agreement.`n"
$ws.Range("C45").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie to be shown on Friday.`n"
$ws.Range("C46").Value = "MSG: None`n`nMSG: The decision has been made to acquire the rights for both movies.`n"
$ws.Range("D46").Value = "both_movies, "
$ws.Range("C47").Value = "MSG: None`n`nMSG: The decision has been recorded as no selection for the movie at this time.`n"
$ws.Range("C48").Value = "MSG: None`n`nMSG: The decision-making process did not result in a choice of movie for Friday.`n"
$ws.Range("C49").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision made regarding Friday's movie.`n"
$ws.Range("C50").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie has been selected for Friday.`n"
$ws.Range("C51").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected in the meeting.`n"
$ws.Range("C52").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected for Friday.`n"
$ws.Range("C53").Value = "MSG: None`n`nMSG: The decision regarding which movie to show on Friday has resulted in no agreement.`n"
$ws.Range("C54").Value = "MSG: None`n`nMSG: The decision has been recorded, and the rights for `"Barbie`" will be acquired.`n"
$ws.Range("C55").Value = "MSG: None`n`nMSG: The rights to both `"Oppenheimer`" and `"Barbie`" have been successfully acquired for the Friday screening.`n"
$ws.Range("D55").Value = "both_movies, "
$ws.Range("C56").Value = "MSG: None`n`nMSG: The decision has been recorded as no agreement was reached regarding which movie to select for Friday.`n"
$ws.Range("C57").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie.`"`n"
$ws.Range("C58").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday's showing.`n"
$ws.Range("C59").Value = "MSG: None`n`nMSG: The decision to acquire the rights for `"Barbie`" has been successfully recorded.`n"
$ws.Range("C60").Value = "MSG: None`n`nMSG: The committee has decided to show `"Barbie`" on Friday.`n"
$ws.Range("C61").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not made.`n"
$ws.Range("C62").Value = "MSG: None`n`nMSG: The decision has been recorded as `"no decision`".`n"
$ws.Range("C63").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie selected for Friday.`n"
$ws.Range("C64").Value = "MSG: None`n`nMSG: The decision process concluded without a definitive selection for a movie to show on Friday. Hence, the status is recorded as `"no decision.`"`n"
$ws.Range("C65").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie was selected during the discussion.`n"
$ws.Range("C66").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie to show on Friday.`n"
$ws.Range("C67").Value = "MSG: None`n`nMSG: The decision regarding the movie to be shown on Friday ended without a clear agreement, resulting in no movie being selected.`n"
$ws.Range("C68").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday has been recorded as having no consensus.`n"
$ws.Range("C69").Value = "MSG: None`n`nMSG: The decision to acquire a movie has not been made, so there will be no selection for Friday.`n"
